$wb = $excel.ActiveWorkbook

# Sheet "展览" — update "想去人数" (interest count) values in column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 642
$ws1.Range("F6").Value = 9741
$ws1.Range("F8").Value = 332
$ws1.Range("F9").Value = 1233
$ws1.Range("F10").Value = 3331
$ws1.Range("F13").Value = 27
$ws1.Range("F16").Value = 515
$ws1.Range("F17").Value = 105
$ws1.Range("F19").Value = 1407

# Sheet "全部类型" — same underlying rows, shifted by +1 due to an extra row
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 642
$ws4.Range("F7").Value = 9741
$ws4.Range("F9").Value = 332
$ws4.Range("F10").Value = 1233
$ws4.Range("F11").Value = 3331
$ws4.Range("F14").Value = 27
$ws4.Range("F17").Value = 515
$ws4.Range("F18").Value = 105
$ws4.Range("F20").Value = 1407
